# Update the workbook metadata / content to reflect the LinuxForHealth
# rebrand (was "Alvearie") and the new 8.0.0 build.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: Property / Value table -----------------------------
$meta = $wb.Worksheets.Item("Metadata")

$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/exchange-plan-indicator"

$meta.Range("B2").Value = $newUrl
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet: the FHIR element table ------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The URL also appears in the "Extension.url" row's Fixed Value column (Q5) -
# keep it in sync with the Metadata sheet's URL above.
$elements.Range("Q5").Value = $newUrl

# The "Extension" row's combined ele-1/ext-1 constraint note is no longer
# shown directly on the base "Extension" row (it now only appears on the
# "Extension.extension" row further down), so clear it here.
$elements.Range("AI2").ClearContents()
